# LOQ4050.xlsx edit
#
# The original row 13 (A13 empty; B13/C13 held
# "5840671 - Francisco José Moreira Chaves") is removed. That shifts every
# row below it up by one, carrying row heights/styles along automatically.
# After the shift a handful of cells need their text corrected to reuse
# strings that already exist elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 13 entirely; rows 14-24 shift up to become rows 13-23.
$ws.Rows.Item(13).Delete()

# Plain text reassignments (no risk of Excel's auto date/number coercion).
$ws.Range("B10").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C10").Value = "5840671 - Francisco José Moreira Chaves"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B18").Value = "5840671 - Francisco José Moreira Chaves"
$ws.Range("C18").Value = "5840671 - Francisco José Moreira Chaves"

$ws.Range("B19").Value = "Por meio de aulas presenciais, com apresentação dos fundamentos, e resolução de exercícios e exemplos aplicativos com uso de tabelas e normas específicas."
$ws.Range("C19").Value = "Por meio de aulas presenciais, com apresentação dos fundamentos, e resolução de exercícios e exemplos aplicativos com uso de tabelas e normas específicas."

$ws.Range("B20").Value = "A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2"
$ws.Range("C20").Value = "A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2"

$ws.Range("B21").Value = "Prova de exame."
$ws.Range("C21").Value = "Prova de exame."

# B15/C15 need to become the literal text "01/01/2011" (same text already
# used in A8/B8/C8). Assigning that string straight into .Value would be
# auto-parsed into a date serial by Excel, so instead copy the already-text
# cell B8 and paste values-only into B15/C15 - this preserves the original
# cell formatting/style while swapping in the text content untouched.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)
